# Revert the SBOL part-URL values in row 3 back to their old order:
#   E3: B0032 -> B0015
#   F3: B0015 -> ComponentDefinition_dvk_backbone_core
#   G3: ComponentDefinition_dvk_backbone_core -> B0032
# (the A5/A6/B1 labels in row 4 stay exactly where they are)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/B0015/1"
$ws.Range("F3").Value = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/ComponentDefinition_dvk_backbone_core/1"
$ws.Range("G3").Value = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/B0032/1"
